# Re-export playlist_per_work/krul001rosi01.xlsx with no is_pref and no lev distance.
#
# The underlying data set was regenerated: the (id, speaker_variant) pairs in
# columns B/C are reshuffled row-for-row, and the "is_prefered" marker ("x")
# that used to flag the first 14 candidate rows (2-15) in column D is dropped
# entirely, since preference is no longer computed for this export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#argire"
$ws.Range("C2").Value = "Argire"
$ws.Range("B3").Value = "#rosanier"
$ws.Range("C3").Value = "Rosanier"
$ws.Range("B4").Value = "#polican:"
$ws.Range("C4").Value = "polican:"
$ws.Range("B5").Value = "#konink"
$ws.Range("C5").Value = "Konink"
$ws.Range("B6").Value = "#cephisi"
$ws.Range("C6").Value = "Cephisi"
$ws.Range("B7").Value = "#bellinde"
$ws.Range("C7").Value = "Bellinde"
$ws.Range("B8").Value = "#chephise"
$ws.Range("C8").Value = "Chephise"
$ws.Range("B9").Value = "#ergun:"
$ws.Range("C9").Value = "Ergun:"
$ws.Range("B10").Value = "#konink:"
$ws.Range("C10").Value = "Konink:"
$ws.Range("B11").Value = "#oront"
$ws.Range("C11").Value = "Oront"
$ws.Range("B12").Value = "#argire:"
$ws.Range("C12").Value = "Argire:"
$ws.Range("B13").Value = "#raetshe"
$ws.Range("C13").Value = "Raetshe"
$ws.Range("B14").Value = "#oronte"
$ws.Range("C14").Value = "Oronte"
$ws.Range("B15").Value = "#oront;"
$ws.Range("C15").Value = "Oront;"
$ws.Range("B16").Value = "#hoveling"
$ws.Range("C16").Value = "Hoveling"
$ws.Range("B17").Value = "#oront:"
$ws.Range("C17").Value = "Oront:"
$ws.Range("B18").Value = "#raedtsh"
$ws.Range("C18").Value = "Raedtsh"
$ws.Range("B19").Value = "#cephise"
$ws.Range("C19").Value = "Cephise"
$ws.Range("B20").Value = "#oront,"
$ws.Range("C20").Value = "Oront,"
$ws.Range("B21").Value = "#verance"
$ws.Range("C21").Value = "Verance"
$ws.Range("B22").Value = "#polican:"
$ws.Range("C22").Value = "Polican:"
$ws.Range("B23").Value = "#rosilion"
$ws.Range("C23").Value = "Rosilion"
$ws.Range("B24").Value = "#rosani"
$ws.Range("C24").Value = "Rosani"
$ws.Range("B25").Value = "#ergunt:"
$ws.Range("C25").Value = "Ergunt:"
$ws.Range("B26").Value = "#cephis"
$ws.Range("C26").Value = "Cephis"
$ws.Range("B27").Value = "#oronte:"
$ws.Range("C27").Value = "Oronte:"
$ws.Range("B28").Value = "#rosil"
$ws.Range("C28").Value = "Rosil"
$ws.Range("B29").Value = "#celiodant"
$ws.Range("C29").Value = "Celiodant"
$ws.Range("B30").Value = "#rosani:"
$ws.Range("C30").Value = "Rosani:"

# Drop the "is_prefered" = "x" marker from D2:D15, leaving each cell present
# but blank (an empty text cell, same shape as the already-blank D16:D30)
# rather than deleting the cell outright.
for ($r = 2; $r -le 15; $r++) {
    $ws.Range("D$r").Formula = "'"
    $ws.Range("D$r").Style = "Normal"
}
